$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "545.48") as well as
# strings with multiple dots (e.g. "60.478.18") that must stay literal text,
# matching the workbook author's original inline-string formatting. Forcing
# the cell to Text format before assigning keeps Excel from auto-converting
# numeric-looking values (dropping trailing zeros, etc.); resetting the style
# back to Normal afterwards avoids leaving a stray custom number format.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.478.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.334.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.331.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.745.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.411.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.327.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0734"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("E34").Value = "  +10.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.381"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "322.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0940"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.564"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0496"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0214"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0213"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.53%  "
